$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2021" column (K), mirroring the formatting
# of the existing "2020" column (J) for each row.
$ws.Range("J3:J10").Copy()
$ws.Range("K3:K10").PasteSpecial(-4122)

# Fill in the new column's values.
$ws.Cells.Item(4, 11).Value = 2021
$ws.Cells.Item(5, 11).Value = 375
$ws.Cells.Item(6, 11).Value = "-"
$ws.Cells.Item(7, 11).Value = 5
$ws.Cells.Item(8, 11).Value = "-"
$ws.Cells.Item(9, 11).Value = 18
$ws.Cells.Item(10, 11).Value = 150

# Match the saved selection in the source workbook.
$ws.Range("K7").Select()
